# Weekly update: insert a new week's worth of price rows (Primera / Segunda)
# at the top of the data block, pushing the existing history down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the first data-shift point (row 296),
# shifting rows 296:435 down to 298:437.
$ws.Rows("296:297").Insert()

# Seed the two new rows with a copy of the week that is now sitting at
# rows 298:299 (i.e. what used to be rows 296:297 before the insert), so all
# the non-date columns carry over unchanged.
$ws.Rows("298:299").Copy()
$ws.Rows("296:297").PasteSpecial()

# Stamp the new rows with the new week's date (2023-04-11).
$ws.Cells.Item(296, 4).Value = 45027
$ws.Cells.Item(297, 4).Value = 45027
